# Applies the row-permutation edit to rows 29-41 of the "Artfynd" sheet.
# Each physical row 29-41 now carries a different field observation's data
# (rows were re-keyed/re-ordered upstream); this script rewrites every cell
# of each row to match. "Startdatum"/"Slutdatum" (Y/AA) are "2023-09-17" for
# every row in this block both before and after, so those two columns are
# intentionally left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the block first (skip Y/AA - unchanged), so stale cells from the
# previous row layout do not linger on columns the new row does not use.
# (Issued as 3 separate single-area calls - ClearContents on a multi-area
# union only reliably clears the first area.)
$ws.Range("A29:X41").ClearContents()
$ws.Range("Z29:Z41").ClearContents()
$ws.Range("AB29:AY41").ClearContents()

# Row 29
$ws.Range("A29").Value = 112170180
$ws.Range("B29").Value = 94173
$ws.Range("C29").Value = "Ovaliderad"
$ws.Range("D29").Value = "LC"
$ws.Range("E29").Value = 2590
$ws.Range("F29").Value = "Kornknutmossa"
$ws.Range("G29").Value = "Odontoschisma denudatum"
$ws.Range("H29").Value = "(Mart.) Dumort"
$ws.Range("P29").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q29").Value = 693927
$ws.Range("R29").Value = 6552111
$ws.Range("S29").Value = 10
$ws.Range("T29").Value = "Stockholm"
$ws.Range("U29").Value = "Haninge"
$ws.Range("V29").Value = "Södermanland"
$ws.Range("W29").Value = "Ornö"
$ws.Range("Z29").Value = "12:10"
$ws.Range("AB29").Value = "12:10"
$ws.Range("AD29").Value = $false
$ws.Range("AE29").Value = $false
$ws.Range("AG29").Value = $false
$ws.Range("AO29").Value = "Kraftigt nedbruten klen låga i sumpskog."
$ws.Range("AW29").Value = "Klas Magnusson"
$ws.Range("AX29").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 30
$ws.Range("A30").Value = 112170185
$ws.Range("B30").Value = 89425
$ws.Range("C30").Value = "Ovaliderad"
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 5442
$ws.Range("F30").Value = "Tallticka"
$ws.Range("G30").Value = "Porodaedalea pini"
$ws.Range("H30").Value = "(Brot.) Murrill"
$ws.Range("P30").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q30").Value = 694000
$ws.Range("R30").Value = 6552145
$ws.Range("S30").Value = 10
$ws.Range("T30").Value = "Stockholm"
$ws.Range("U30").Value = "Haninge"
$ws.Range("V30").Value = "Södermanland"
$ws.Range("W30").Value = "Ornö"
$ws.Range("Z30").Value = "11:45"
$ws.Range("AB30").Value = "11:45"
$ws.Range("AD30").Value = $false
$ws.Range("AE30").Value = $false
$ws.Range("AG30").Value = $false
$ws.Range("AI30").Value = "Tallskog"
$ws.Range("AJ30").Value = "tall"
$ws.Range("AK30").Value = "Pinus sylvestris"
$ws.Range("AL30").Value = "Gammal levande tall"
$ws.Range("AO30").Value = "Pinus sylvestris # Gammal levande tall"
$ws.Range("AW30").Value = "Klas Magnusson"
$ws.Range("AX30").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 31
$ws.Range("A31").Value = 112170177
$ws.Range("B31").Value = 93158
$ws.Range("C31").Value = "Ovaliderad"
$ws.Range("D31").Value = "LC"
$ws.Range("E31").Value = 2667
$ws.Range("F31").Value = "Platt fjädermossa"
$ws.Range("G31").Value = "Neckera complanata"
$ws.Range("H31").Value = "(Hedw.) Huebener"
$ws.Range("P31").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q31").Value = 693885
$ws.Range("R31").Value = 6552076
$ws.Range("S31").Value = 10
$ws.Range("T31").Value = "Stockholm"
$ws.Range("U31").Value = "Haninge"
$ws.Range("V31").Value = "Södermanland"
$ws.Range("W31").Value = "Ornö"
$ws.Range("Z31").Value = "12:29"
$ws.Range("AB31").Value = "12:29"
$ws.Range("AD31").Value = $false
$ws.Range("AE31").Value = $false
$ws.Range("AG31").Value = $false
$ws.Range("AO31").Value = "Vid basen av klippvägg."
$ws.Range("AW31").Value = "Klas Magnusson"
$ws.Range("AX31").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 32
$ws.Range("A32").Value = 112170184
$ws.Range("B32").Value = 93388
$ws.Range("C32").Value = "Ovaliderad"
$ws.Range("D32").Value = "LC"
$ws.Range("E32").Value = 2180
$ws.Range("F32").Value = "Blåmossa"
$ws.Range("G32").Value = "Leucobryum glaucum"
$ws.Range("H32").Value = "(Hedw.) Ångstr."
$ws.Range("P32").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q32").Value = 693996
$ws.Range("R32").Value = 6552132
$ws.Range("S32").Value = 10
$ws.Range("T32").Value = "Stockholm"
$ws.Range("U32").Value = "Haninge"
$ws.Range("V32").Value = "Södermanland"
$ws.Range("W32").Value = "Ornö"
$ws.Range("Z32").Value = "11:50"
$ws.Range("AB32").Value = "11:50"
$ws.Range("AD32").Value = $false
$ws.Range("AE32").Value = $false
$ws.Range("AG32").Value = $false
$ws.Range("AI32").Value = "Tallskog"
$ws.Range("AW32").Value = "Klas Magnusson"
$ws.Range("AX32").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 33
$ws.Range("A33").Value = 112170187
$ws.Range("B33").Value = 93388
$ws.Range("C33").Value = "Ovaliderad"
$ws.Range("D33").Value = "LC"
$ws.Range("E33").Value = 2180
$ws.Range("F33").Value = "Blåmossa"
$ws.Range("G33").Value = "Leucobryum glaucum"
$ws.Range("H33").Value = "(Hedw.) Ångstr."
$ws.Range("P33").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q33").Value = 694042
$ws.Range("R33").Value = 6552162
$ws.Range("S33").Value = 10
$ws.Range("T33").Value = "Stockholm"
$ws.Range("U33").Value = "Haninge"
$ws.Range("V33").Value = "Södermanland"
$ws.Range("W33").Value = "Ornö"
$ws.Range("Z33").Value = "11:30"
$ws.Range("AB33").Value = "11:30"
$ws.Range("AC33").Value = "Stor kudde."
$ws.Range("AD33").Value = $false
$ws.Range("AE33").Value = $false
$ws.Range("AG33").Value = $false
$ws.Range("AW33").Value = "Klas Magnusson"
$ws.Range("AX33").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 34
$ws.Range("A34").Value = 112170186
$ws.Range("B34").Value = 78107
$ws.Range("C34").Value = "Ovaliderad"
$ws.Range("D34").Value = "NT"
$ws.Range("E34").Value = 6453
$ws.Range("F34").Value = "Vedskivlav"
$ws.Range("G34").Value = "Hertelidea botryosa"
$ws.Range("H34").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("P34").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q34").Value = 694041
$ws.Range("R34").Value = 6552157
$ws.Range("S34").Value = 10
$ws.Range("T34").Value = "Stockholm"
$ws.Range("U34").Value = "Haninge"
$ws.Range("V34").Value = "Södermanland"
$ws.Range("W34").Value = "Ornö"
$ws.Range("Z34").Value = "11:37"
$ws.Range("AB34").Value = "11:37"
$ws.Range("AD34").Value = $false
$ws.Range("AE34").Value = $false
$ws.Range("AG34").Value = $false
$ws.Range("AI34").Value = "Hällmarkstallskog"
$ws.Range("AJ34").Value = "tall"
$ws.Range("AK34").Value = "Pinus sylvestris"
$ws.Range("AL34").Value = "Gammal torr tallåga"
$ws.Range("AO34").Value = "Pinus sylvestris # Gammal torr tallåga"
$ws.Range("AW34").Value = "Klas Magnusson"
$ws.Range("AX34").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 35
$ws.Range("A35").Value = 112170182
$ws.Range("B35").Value = 8367
$ws.Range("C35").Value = "Ovaliderad"
$ws.Range("D35").Value = "LC"
$ws.Range("E35").Value = 106554
$ws.Range("F35").Value = "Björksplintborre"
$ws.Range("G35").Value = "Scolytus ratzeburgii"
$ws.Range("H35").Value = "Janson, 1856"
$ws.Range("M35").Value = "äldre gnagspår"
$ws.Range("P35").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q35").Value = 693976
$ws.Range("R35").Value = 6552129
$ws.Range("S35").Value = 10
$ws.Range("T35").Value = "Stockholm"
$ws.Range("U35").Value = "Haninge"
$ws.Range("V35").Value = "Södermanland"
$ws.Range("W35").Value = "Ornö"
$ws.Range("Z35").Value = "11:55"
$ws.Range("AB35").Value = "11:55"
$ws.Range("AD35").Value = $false
$ws.Range("AE35").Value = $false
$ws.Range("AG35").Value = $false
$ws.Range("AJ35").Value = "björkar"
$ws.Range("AK35").Value = "Betula"
$ws.Range("AL35").Value = "Björklåga med det mesta av barken kvar."
$ws.Range("AO35").Value = "Betula # Björklåga med det mesta av barken kvar."
$ws.Range("AW35").Value = "Klas Magnusson"
$ws.Range("AX35").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 36
$ws.Range("A36").Value = 112170178
$ws.Range("B36").Value = 93171
$ws.Range("C36").Value = "Ovaliderad"
$ws.Range("D36").Value = "LC"
$ws.Range("E36").Value = 2818
$ws.Range("F36").Value = "Stubbspretmossa"
$ws.Range("G36").Value = "Herzogiella seligeri"
$ws.Range("H36").Value = "(Brid.) Z.Iwats."
$ws.Range("P36").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q36").Value = 693919
$ws.Range("R36").Value = 6552120
$ws.Range("S36").Value = 10
$ws.Range("T36").Value = "Stockholm"
$ws.Range("U36").Value = "Haninge"
$ws.Range("V36").Value = "Södermanland"
$ws.Range("W36").Value = "Ornö"
$ws.Range("Z36").Value = "12:19"
$ws.Range("AB36").Value = "12:19"
$ws.Range("AD36").Value = $false
$ws.Range("AE36").Value = $false
$ws.Range("AG36").Value = $false
$ws.Range("AW36").Value = "Klas Magnusson"
$ws.Range("AX36").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 37
$ws.Range("A37").Value = 112170183
$ws.Range("B37").Value = 93388
$ws.Range("C37").Value = "Ovaliderad"
$ws.Range("D37").Value = "LC"
$ws.Range("E37").Value = 2180
$ws.Range("F37").Value = "Blåmossa"
$ws.Range("G37").Value = "Leucobryum glaucum"
$ws.Range("H37").Value = "(Hedw.) Ångstr."
$ws.Range("P37").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q37").Value = 693999
$ws.Range("R37").Value = 6552113
$ws.Range("S37").Value = 10
$ws.Range("T37").Value = "Stockholm"
$ws.Range("U37").Value = "Haninge"
$ws.Range("V37").Value = "Södermanland"
$ws.Range("W37").Value = "Ornö"
$ws.Range("Z37").Value = "11:53"
$ws.Range("AB37").Value = "11:53"
$ws.Range("AC37").Value = "Flera stora kuddar i blåbärsris."
$ws.Range("AD37").Value = $false
$ws.Range("AE37").Value = $false
$ws.Range("AG37").Value = $false
$ws.Range("AI37").Value = "Tallskog"
$ws.Range("AW37").Value = "Klas Magnusson"
$ws.Range("AX37").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 38
$ws.Range("A38").Value = 112170190
$ws.Range("B38").Value = 5113
$ws.Range("C38").Value = "Ovaliderad"
$ws.Range("D38").Value = "LC"
$ws.Range("E38").Value = 100526
$ws.Range("F38").Value = "Bronshjon"
$ws.Range("G38").Value = "Callidium coriaceum"
$ws.Range("H38").Value = "Paykull, 1800"
$ws.Range("M38").Value = "äldre gnagspår"
$ws.Range("P38").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q38").Value = 694070
$ws.Range("R38").Value = 6552221
$ws.Range("S38").Value = 10
$ws.Range("T38").Value = "Stockholm"
$ws.Range("U38").Value = "Haninge"
$ws.Range("V38").Value = "Södermanland"
$ws.Range("W38").Value = "Ornö"
$ws.Range("Z38").Value = "10:50"
$ws.Range("AB38").Value = "10:50"
$ws.Range("AD38").Value = $false
$ws.Range("AE38").Value = $false
$ws.Range("AG38").Value = $false
$ws.Range("AJ38").Value = "gran"
$ws.Range("AK38").Value = "Picea abies"
$ws.Range("AL38").Value = "Granlåga"
$ws.Range("AO38").Value = "Picea abies # Granlåga"
$ws.Range("AW38").Value = "Klas Magnusson"
$ws.Range("AX38").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 39
$ws.Range("A39").Value = 112170188
$ws.Range("B39").Value = 4717
$ws.Range("C39").Value = "Ovaliderad"
$ws.Range("D39").Value = "LC"
$ws.Range("E39").Value = 102306
$ws.Range("F39").Value = "Granbarkgnagare"
$ws.Range("G39").Value = "Microbregma emarginatum"
$ws.Range("H39").Value = "(Duftschmid, 1825)"
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = "1"
$ws.Range("J39").Value = "ex."
$ws.Range("K39").Value = "larv/nymf"
$ws.Range("P39").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q39").Value = 694095
$ws.Range("R39").Value = 6552131
$ws.Range("S39").Value = 10
$ws.Range("T39").Value = "Stockholm"
$ws.Range("U39").Value = "Haninge"
$ws.Range("V39").Value = "Södermanland"
$ws.Range("W39").Value = "Ornö"
$ws.Range("Z39").Value = "11:06"
$ws.Range("AB39").Value = "11:06"
$ws.Range("AD39").Value = $false
$ws.Range("AE39").Value = $false
$ws.Range("AG39").Value = $false
$ws.Range("AJ39").Value = "gran"
$ws.Range("AK39").Value = "Picea abies"
$ws.Range("AL39").Value = "Äldre grovbarkig levande gran."
$ws.Range("AO39").Value = "Picea abies # Äldre grovbarkig levande gran."
$ws.Range("AW39").Value = "Klas Magnusson"
$ws.Range("AX39").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 40
$ws.Range("A40").Value = 112170151
$ws.Range("B40").Value = 56414
$ws.Range("C40").Value = "Ovaliderad"
$ws.Range("D40").Value = "NT"
$ws.Range("E40").Value = 100049
$ws.Range("F40").Value = "Spillkråka"
$ws.Range("G40").Value = "Dryocopus martius"
$ws.Range("H40").Value = "(Linnaeus, 1758)"
$ws.Range("P40").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q40").Value = 693854
$ws.Range("R40").Value = 6552210
$ws.Range("S40").Value = 75
$ws.Range("T40").Value = "Stockholm"
$ws.Range("U40").Value = "Haninge"
$ws.Range("V40").Value = "Södermanland"
$ws.Range("W40").Value = "Ornö"
$ws.Range("Z40").Value = "16:11"
$ws.Range("AB40").Value = "16:11"
$ws.Range("AC40").Value = "Läte nere vid vattnet."
$ws.Range("AD40").Value = $false
$ws.Range("AE40").Value = $false
$ws.Range("AG40").Value = $false
$ws.Range("AW40").Value = "Klas Magnusson"
$ws.Range("AX40").Value = "Klas Magnusson, Måns Persson, Per Flodby"

# Row 41
$ws.Range("A41").Value = 112170191
$ws.Range("B41").Value = 93388
$ws.Range("C41").Value = "Ovaliderad"
$ws.Range("D41").Value = "LC"
$ws.Range("E41").Value = 2180
$ws.Range("F41").Value = "Blåmossa"
$ws.Range("G41").Value = "Leucobryum glaucum"
$ws.Range("H41").Value = "(Hedw.) Ångstr."
$ws.Range("P41").Value = "Gråberget, Ornö, Srm"
$ws.Range("Q41").Value = 694067
$ws.Range("R41").Value = 6552283
$ws.Range("S41").Value = 10
$ws.Range("T41").Value = "Stockholm"
$ws.Range("U41").Value = "Haninge"
$ws.Range("V41").Value = "Södermanland"
$ws.Range("W41").Value = "Ornö"
$ws.Range("Z41").Value = "10:36"
$ws.Range("AB41").Value = "10:36"
$ws.Range("AC41").Value = "Stor kudde."
$ws.Range("AD41").Value = $false
$ws.Range("AE41").Value = $false
$ws.Range("AG41").Value = $false
$ws.Range("AI41").Value = "Fuktig granskog med tall och björk."
$ws.Range("AW41").Value = "Klas Magnusson"
$ws.Range("AX41").Value = "Klas Magnusson, Måns Persson, Per Flodby"
